$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 347.57144
$ws.Range("J2").Value = 499
$ws.Range("L2").Value = 499
$ws.Range("N2").Value = -725

$ws.Range("H28").Value = 1152.2162
$ws.Range("I28").Value = 681.70966
$ws.Range("K28").Value = 681.70966
$ws.Range("M28").Value = -196.70966

$ws.Range("H40").Value = 4205.2085
$ws.Range("I40").Value = 3799.913
$ws.Range("J40").Value = 4578.08
$ws.Range("K40").Value = 3799.913
$ws.Range("L40").Value = 4578.08
$ws.Range("M40").Value = -3624.913
$ws.Range("N40").Value = -4928.08

$ws.Range("H62").Value = 5846.5264
$ws.Range("I62").Value = 2868.4285
$ws.Range("K62").Value = 2868.4285
$ws.Range("M62").Value = -2244.4285

$ws.Range("H65").Value = 5846.5264
$ws.Range("I65").Value = 2868.4285
$ws.Range("K65").Value = 14342.1425
$ws.Range("M65").Value = -11222.1425

$ws.Range("H70").Value = 29301.1
$ws.Range("I70").Value = 104388.8
$ws.Range("J70").Value = 4271.8667
$ws.Range("K70").Value = 313166.4
$ws.Range("L70").Value = 12815.6001
$ws.Range("M70").Value = -312896.4
$ws.Range("N70").Value = -13355.6001

$ws.Range("H73").Value = 29301.1
$ws.Range("I73").Value = 104388.8
$ws.Range("J73").Value = 4271.8667
$ws.Range("K73").Value = 313166.4
$ws.Range("L73").Value = 12815.6001
$ws.Range("M73").Value = -312230.4
$ws.Range("N73").Value = -14687.6001

$ws.Range("H74").Value = 6508.864
$ws.Range("I74").Value = 5168.091
$ws.Range("J74").Value = 7849.636
$ws.Range("K74").Value = 5168.091
$ws.Range("L74").Value = 7849.636
$ws.Range("M74").Value = -4232.091
$ws.Range("N74").Value = -9721.636

$ws.Range("H77").Value = 6508.864
$ws.Range("I77").Value = 5168.091
$ws.Range("J77").Value = 7849.636
$ws.Range("K77").Value = 25840.455
$ws.Range("L77").Value = 39248.18
$ws.Range("M77").Value = -21160.455
$ws.Range("N77").Value = -48608.18

$ws.Range("H80").Value = 1803.1818
$ws.Range("I80").Value = 550.8333
$ws.Range("J80").Value = 2272.8125
$ws.Range("K80").Value = 1652.4999
$ws.Range("L80").Value = 6818.4375
$ws.Range("M80").Value = -654.4999
$ws.Range("N80").Value = -8814.4375

$ws.Range("H83").Value = 1803.1818
$ws.Range("I83").Value = 550.8333
$ws.Range("J83").Value = 2272.8125
$ws.Range("K83").Value = 4957.4997
$ws.Range("L83").Value = 20455.3125
$ws.Range("M83").Value = 34.5002999999997
$ws.Range("N83").Value = -30439.3125

$ws.Range("H92").Value = 2814.5625
$ws.Range("I92").Value = 1512.2727
$ws.Range("K92").Value = 1512.2727
$ws.Range("M92").Value = -264.2727

$ws.Range("H98").Value = 383162.62
$ws.Range("I98").Value = 916.7273
$ws.Range("K98").Value = 916.7273
$ws.Range("M98").Value = 581.2727

$ws.Range("H116").Value = 12653.23
$ws.Range("I116").Value = 12387.667
$ws.Range("K116").Value = 12387.667
$ws.Range("M116").Value = -8945.666999999999

$ws.Range("H122").Value = 383162.62
$ws.Range("I122").Value = 916.7273
$ws.Range("K122").Value = 2750.1819
$ws.Range("M122").Value = -300.1819

$ws.Range("H131").Value = 2454.1538
$ws.Range("I131").Value = 671.7895
$ws.Range("K131").Value = 2015.3685
$ws.Range("M131").Value = 3024.6315

$ws.Range("H132").Value = 1851.8636
$ws.Range("I132").Value = 1811
$ws.Range("J132").Value = 2110.6667
$ws.Range("K132").Value = 5433
$ws.Range("L132").Value = 6332.000100000001
$ws.Range("M132").Value = -2903
$ws.Range("N132").Value = -11392.0001

$ws.Range("H137").Value = 17243800
$ws.Range("I137").Value = 50002020
$ws.Range("J137").Value = 2631.5789
$ws.Range("K137").Value = 150006060
$ws.Range("L137").Value = 7894.736699999999
$ws.Range("M137").Value = -150003510
$ws.Range("N137").Value = -12994.7367

$ws.Range("H138").Value = 3650.7307
$ws.Range("I138").Value = 2413.3076
$ws.Range("J138").Value = 4063.205
$ws.Range("K138").Value = 7239.9228
$ws.Range("L138").Value = 12189.615
$ws.Range("M138").Value = -2099.9228
$ws.Range("N138").Value = -22469.615

$ws.Range("H141").Value = 1326.3889
$ws.Range("J141").Value = 3998.6667
$ws.Range("L141").Value = 11996.0001
$ws.Range("N141").Value = -22356.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3354.2183
$ws.Range("I32").Value = 3354.2183
$ws.Range("K32").Value = 3354.2183
$ws.Range("M32").Value = -3067.2183

$ws.Range("H61").Value = 10377.75
$ws.Range("I61").Value = 3848.111
$ws.Range("J61").Value = 18773
$ws.Range("K61").Value = 3848.111
$ws.Range("L61").Value = 18773
$ws.Range("M61").Value = -3636.111
$ws.Range("N61").Value = -19197

$ws.Range("H74").Value = 8774740
$ws.Range("I74").Value = 9805879
$ws.Range("J74").Value = 10065.75
$ws.Range("K74").Value = 9805879
$ws.Range("L74").Value = 10065.75
$ws.Range("M74").Value = -9805005
$ws.Range("N74").Value = -11813.75

$ws.Range("H77").Value = 8774740
$ws.Range("I77").Value = 9805879
$ws.Range("J77").Value = 10065.75
$ws.Range("K77").Value = 49029395
$ws.Range("L77").Value = 50328.75
$ws.Range("M77").Value = -49025027
$ws.Range("N77").Value = -59064.75

$ws.Range("H110").Value = 4167.9287
$ws.Range("I110").Value = 1851.8
$ws.Range("K110").Value = 1851.8
$ws.Range("M110").Value = 193.2

$ws.Range("H122").Value = 55559556
$ws.Range("I122").Value = 2355.8572
$ws.Range("J122").Value = 90914136
$ws.Range("K122").Value = 7067.571599999999
$ws.Range("L122").Value = 272742408
$ws.Range("M122").Value = -4617.571599999999
$ws.Range("N122").Value = -272747308

$ws.Range("H125").Value = 447800
$ws.Range("J125").Value = 447800
$ws.Range("L125").Value = 447800
$ws.Range("N125").Value = -457640

$ws.Range("H132").Value = 4323.15
$ws.Range("I132").Value = 3549.9736
$ws.Range("J132").Value = 19013.5
$ws.Range("K132").Value = 10649.9208
$ws.Range("L132").Value = 57040.5
$ws.Range("M132").Value = -8119.9208
$ws.Range("N132").Value = -62100.5

$ws.Range("H136").Value = 10377.75
$ws.Range("I136").Value = 3848.111
$ws.Range("J136").Value = 18773
$ws.Range("K136").Value = 11544.333
$ws.Range("L136").Value = 56319
$ws.Range("M136").Value = -8994.332999999999
$ws.Range("N136").Value = -61419

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 727.375
$ws.Range("I22").Value = 738.4286
$ws.Range("J22").Value = 650
$ws.Range("K22").Value = 738.4286
$ws.Range("L22").Value = 650
$ws.Range("M22").Value = -565.4286
$ws.Range("N22").Value = -996

$ws.Range("H80").Value = 497.6316
$ws.Range("J80").Value = 484
$ws.Range("L80").Value = 484
$ws.Range("N80").Value = -2480

$ws.Range("H83").Value = 497.6316
$ws.Range("J83").Value = 484
$ws.Range("L83").Value = 2420
$ws.Range("N83").Value = -12404

$ws.Range("H86").Value = 2570.353
$ws.Range("I86").Value = 2168.1155
$ws.Range("J86").Value = 3877.625
$ws.Range("K86").Value = 2168.1155
$ws.Range("L86").Value = 3877.625
$ws.Range("M86").Value = -1045.1155
$ws.Range("N86").Value = -6123.625

$ws.Range("H89").Value = 2570.353
$ws.Range("I89").Value = 2168.1155
$ws.Range("J89").Value = 3877.625
$ws.Range("K89").Value = 10840.5775
$ws.Range("L89").Value = 19388.125
$ws.Range("M89").Value = -5224.577499999999
$ws.Range("N89").Value = -30620.125

$ws.Range("H105").Value = 23120.385
$ws.Range("I105").Value = 30260.143
$ws.Range("K105").Value = 30260.143
$ws.Range("M105").Value = -28513.143

$ws.Range("H111").Value = 50000
$ws.Range("J111").Value = 50000
$ws.Range("L111").Value = 50000
$ws.Range("N111").Value = -58180

$ws.Range("H134").Value = 1521.25
$ws.Range("I134").Value = 1025.2941
$ws.Range("J134").Value = 4331.6665
$ws.Range("K134").Value = 3075.8823
$ws.Range("L134").Value = 12994.9995
$ws.Range("M134").Value = -540.8823000000002
$ws.Range("N134").Value = -18064.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1777.375
$ws.Range("I22").Value = 260
$ws.Range("J22").Value = 3294.75
$ws.Range("K22").Value = 260
$ws.Range("L22").Value = 3294.75
$ws.Range("M22").Value = 90
$ws.Range("N22").Value = -3994.75

$ws.Range("H23").Value = 7498.6113
$ws.Range("I23").Value = 5165
$ws.Range("J23").Value = 19166.666
$ws.Range("K23").Value = 5165
$ws.Range("L23").Value = 19166.666
$ws.Range("M23").Value = -4925
$ws.Range("N23").Value = -19646.666

$ws.Range("H27").Value = 7498.6113
$ws.Range("I27").Value = 5165
$ws.Range("J27").Value = 19166.666
$ws.Range("K27").Value = 5165
$ws.Range("L27").Value = 19166.666
$ws.Range("M27").Value = -4973
$ws.Range("N27").Value = -19550.666

$ws.Range("H31").Value = 31424.256
$ws.Range("I31").Value = 3227.6924
$ws.Range("J31").Value = 87817.38
$ws.Range("K31").Value = 3227.6924
$ws.Range("L31").Value = 87817.38
$ws.Range("M31").Value = -2932.6924
$ws.Range("N31").Value = -88407.38

$ws.Range("H34").Value = 31424.256
$ws.Range("I34").Value = 3227.6924
$ws.Range("J34").Value = 87817.38
$ws.Range("K34").Value = 3227.6924
$ws.Range("L34").Value = 87817.38
$ws.Range("M34").Value = -3025.6924
$ws.Range("N34").Value = -88221.38

$ws.Range("H58").Value = 4173.778
$ws.Range("I58").Value = 2128.45
$ws.Range("K58").Value = 2128.45
$ws.Range("M58").Value = -1925.45

$ws.Range("H86").Value = 6749.75
$ws.Range("I86").Value = 6499.5
$ws.Range("K86").Value = 6499.5
$ws.Range("M86").Value = -5376.5

$ws.Range("H89").Value = 6749.75
$ws.Range("I89").Value = 6499.5
$ws.Range("K89").Value = 32497.5
$ws.Range("M89").Value = -26881.5

$ws.Range("H99").Value = 5187
$ws.Range("I99").Value = 3966.25
$ws.Range("J99").Value = 7628.5
$ws.Range("K99").Value = 3966.25
$ws.Range("L99").Value = 7628.5
$ws.Range("M99").Value = -2468.25
$ws.Range("N99").Value = -10624.5

$ws.Range("H107").Value = 1966
$ws.Range("I107").Value = 1634.3334
$ws.Range("K107").Value = 1634.3334
$ws.Range("M107").Value = 285.6666

$ws.Range("H126").Value = 5187
$ws.Range("I126").Value = 3966.25
$ws.Range("J126").Value = 7628.5
$ws.Range("K126").Value = 11898.75
$ws.Range("L126").Value = 22885.5
$ws.Range("M126").Value = -9428.75
$ws.Range("N126").Value = -27825.5

$ws.Range("H132").Value = 3979.3572
$ws.Range("I132").Value = 2724.75
$ws.Range("J132").Value = 11507
$ws.Range("K132").Value = 8174.25
$ws.Range("L132").Value = 34521
$ws.Range("M132").Value = -5644.25
$ws.Range("N132").Value = -39581

$ws.Range("H134").Value = 7460.405
$ws.Range("I134").Value = 7054.1377
$ws.Range("J134").Value = 8366.691999999999
$ws.Range("K134").Value = 21162.4131
$ws.Range("L134").Value = 25100.076
$ws.Range("M134").Value = -18627.4131
$ws.Range("N134").Value = -30170.076

$ws.Range("H136").Value = 4173.778
$ws.Range("I136").Value = 2128.45
$ws.Range("K136").Value = 6385.349999999999
$ws.Range("M136").Value = -3835.349999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 4802
$ws.Range("I58").Value = 4502.5
$ws.Range("K58").Value = 13507.5
$ws.Range("M58").Value = -13379.5

$ws.Range("H92").Value = 5000.1665
$ws.Range("I92").Value = 2333.3333
$ws.Range("J92").Value = 7667
$ws.Range("K92").Value = 6999.999899999999
$ws.Range("L92").Value = 23001
$ws.Range("M92").Value = -5751.999899999999
$ws.Range("N92").Value = -25497

$ws.Range("H98").Value = 934.2
$ws.Range("I98").Value = 442.14285
$ws.Range("K98").Value = 1326.42855
$ws.Range("M98").Value = 171.5714499999999

$ws.Range("H113").Value = 45455412
$ws.Range("J113").Value = 55556332
$ws.Range("L113").Value = 166668996
$ws.Range("N113").Value = -166673336

$ws.Range("H128").Value = 1058250.2
$ws.Range("I128").Value = 1058250.2
$ws.Range("K128").Value = 3174750.6
$ws.Range("M128").Value = -3169770.6

$ws.Range("H136").Value = 18519752
$ws.Range("I136").Value = 19608852
$ws.Range("J136").Value = 5050
$ws.Range("K136").Value = 58826556
$ws.Range("L136").Value = 15150
$ws.Range("M136").Value = -58821456
$ws.Range("N136").Value = -25350

$ws.Range("H140").Value = 2676.3572
$ws.Range("I140").Value = 2156.25
$ws.Range("K140").Value = 6468.75
$ws.Range("M140").Value = -1288.75

$ws.Range("H141").Value = 7211.5557
$ws.Range("I141").Value = 4135.9443
$ws.Range("J141").Value = 13362.777
$ws.Range("K141").Value = 12407.8329
$ws.Range("L141").Value = 40088.331
$ws.Range("M141").Value = -7227.832900000001
$ws.Range("N141").Value = -50448.331

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7748.75
$ws.Range("I80").Value = 3994.5
$ws.Range("K80").Value = 3994.5
$ws.Range("M80").Value = -2996.5

$ws.Range("H83").Value = 7748.75
$ws.Range("I83").Value = 3994.5
$ws.Range("K83").Value = 19972.5
$ws.Range("M83").Value = -14980.5

$ws.Range("H97").Value = 1304.5862
$ws.Range("J97").Value = 1343
$ws.Range("L97").Value = 1343
$ws.Range("N97").Value = -2335

$ws.Range("H107").Value = 1123.8235
$ws.Range("I107").Value = 593.7143
$ws.Range("K107").Value = 593.7143
$ws.Range("M107").Value = 1326.2857

$ws.Range("H123").Value = 38350.668
$ws.Range("J123").Value = 38350.668
$ws.Range("L123").Value = 38350.668
$ws.Range("N123").Value = -43250.668

$ws.Range("H126").Value = 4997.364
$ws.Range("I126").Value = 3451.6
$ws.Range("K126").Value = 10354.8
$ws.Range("M126").Value = -7884.799999999999

$ws.Range("H132").Value = 4788.7036
$ws.Range("I132").Value = 3558.2727
$ws.Range("J132").Value = 10202.6
$ws.Range("K132").Value = 10674.8181
$ws.Range("L132").Value = 30607.8
$ws.Range("M132").Value = -8144.8181
$ws.Range("N132").Value = -35667.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1069.9286
$ws.Range("I16").Value = 998.44446
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 998.44446
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -828.44446
$ws.Range("N16").Value = -3340

$ws.Range("H21").Value = 2996.5
$ws.Range("J21").Value = 2996.5
$ws.Range("L21").Value = 2996.5
$ws.Range("N21").Value = -3344.5

$ws.Range("H22").Value = 4555.273
$ws.Range("I22").Value = 980.8
$ws.Range("K22").Value = 980.8
$ws.Range("M22").Value = -685.8

$ws.Range("H27").Value = 4555.273
$ws.Range("I27").Value = 980.8
$ws.Range("K27").Value = 980.8
$ws.Range("M27").Value = -873.8

$ws.Range("H40").Value = 5850.0225
$ws.Range("I40").Value = 5769.39
$ws.Range("K40").Value = 5769.39
$ws.Range("M40").Value = -5633.39

$ws.Range("H55").Value = 2383042.5
$ws.Range("I55").Value = 4167215.8
$ws.Range("J55").Value = 4144.5557
$ws.Range("K55").Value = 4167215.8
$ws.Range("L55").Value = 4144.5557
$ws.Range("M55").Value = -4167042.8
$ws.Range("N55").Value = -4490.5557

$ws.Range("H94").Value = 20000
$ws.Range("J94").Value = 20000
$ws.Range("L94").Value = 20000
$ws.Range("N94").Value = -21352

$ws.Range("H122").Value = 153198.92
$ws.Range("I122").Value = 194932.42
$ws.Range("J122").Value = 7131.6665
$ws.Range("K122").Value = 584797.26
$ws.Range("L122").Value = 21394.9995
$ws.Range("M122").Value = -582347.26
$ws.Range("N122").Value = -26294.9995

$ws.Range("H132").Value = 6881.522
$ws.Range("I132").Value = 5781.5293
$ws.Range("J132").Value = 9998.166999999999
$ws.Range("K132").Value = 17344.5879
$ws.Range("L132").Value = 29994.501
$ws.Range("M132").Value = -14814.5879
$ws.Range("N132").Value = -35054.501

$ws.Range("H136").Value = 5106.224
$ws.Range("I136").Value = 3801.4792
$ws.Range("K136").Value = 11404.4376
$ws.Range("M136").Value = -8854.437600000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 2000
$ws.Range("J14").Value = 2000
$ws.Range("L14").Value = 2000
$ws.Range("N14").Value = -2336

$ws.Range("H93").Value = 55495.5
$ws.Range("J93").Value = 55495.5
$ws.Range("L93").Value = 55495.5
$ws.Range("N93").Value = -60487.5

$ws.Range("H100").Value = 349.6154
$ws.Range("J100").Value = 333
$ws.Range("L100").Value = 666
$ws.Range("N100").Value = -1748

$ws.Range("H126").Value = 2929.3333
$ws.Range("I126").Value = 2591.6
$ws.Range("J126").Value = 3170.5715
$ws.Range("K126").Value = 7774.799999999999
$ws.Range("L126").Value = 9511.7145
$ws.Range("M126").Value = -5304.799999999999
$ws.Range("N126").Value = -14451.7145

$ws.Range("H132").Value = 3068.9119
$ws.Range("I132").Value = 1522.2174
$ws.Range("J132").Value = 6302.909
$ws.Range("K132").Value = 4566.6522
$ws.Range("L132").Value = 18908.727
$ws.Range("M132").Value = -2036.6522
$ws.Range("N132").Value = -23968.727

$ws.Range("H136").Value = 2118.4783
$ws.Range("I136").Value = 1376.7561
$ws.Range("K136").Value = 4130.2683
$ws.Range("M136").Value = -1580.2683
